$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (D) and Volume(1h) (E) column updates ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.267.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.02%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.557.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.51%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.555.36"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.496"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.53%  "

$ws.Range("E10").Value = "  -1.59%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.411"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.159.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.18%  "

$ws.Range("E14").Value = "  -2.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.82%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.553.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.76%  "

$ws.Range("E17").Value = "  +2.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.231.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "420.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.55%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.695.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.20%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("E27").Value = "  -1.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.34%  "

$ws.Range("E31").Value = "  +0.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.553.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.10%  "

$ws.Range("E33").Value = "  +2.09%  "

$ws.Range("E34").Value = "  -3.08%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("E36").Value = "  -2.45%  "

$ws.Range("E37").Value = "  -8.55%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.66%  "

$ws.Range("E39").Value = "  -6.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "174.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0828"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.58%  "

$ws.Range("E42").Value = "  -1.38%  "

$ws.Range("E43").Value = "  -2.83%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.74"
$ws.Range("D44").Style = "Normal"

$ws.Range("E45").Value = "  -5.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.54%  "

$ws.Range("E48").Value = "  +0.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.76%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.53%  "

# --- Row 29/30: PancakeSwap and RenderToken swapped position with refreshed values ---
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.58%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.06%  "

Write-Output "Applied cryptos list update"